# Correct the author credit line on every slide:
#   " Bence Mitlasoczki & Lena Gschossmann"
#   -> " Bence Mitlasoczki & Lena J. Gschossmann"
#
# The credit line lives inside a single <a:r> run (together with its
# leading space) in the "Text Placeholder 2" shape's text frame, right
# after the date + en-dash runs. We locate that run's text dynamically
# (rather than hard-coding slide/shape numbers) and only rewrite the
# matched characters, so the surrounding runs (date digits, en-dash,
# spacing, run-level formatting/lang) are left untouched.

$p = $ppt.ActivePresentation

$oldText = " Bence Mitlasoczki & Lena Gschossmann"
$newText = " Bence Mitlasoczki & Lena J. Gschossmann"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if ($shape.HasTextFrame) {
            $textFrame = $shape.TextFrame
            if ($textFrame.HasText) {
                $textRange = $textFrame.TextRange
                $fullText = $textRange.Text

                $matchIndex = $fullText.IndexOf($oldText)
                if ($matchIndex -ge 0) {
                    # Characters() uses 1-based start positions.
                    $startPos = $matchIndex + 1
                    $runRange = $textRange.Characters($startPos, $oldText.Length)
                    $runRange.Text = $newText
                }
            }
        }
    }
}
